# Insert a new weekly data row for "Choclo" (Provincia de Limarí, $/unidad)
# above the old row 780. Excel's native row-insert shifts every row from
# 780 downward to 781..821 (and extends the sheet dimension to R821),
# exactly like the diff shows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("780:780").Insert()

$ws.Range("A780").Value = 8
$ws.Range("B780").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C780").Value = 'Coquimbo'
$ws.Range("D780").Value = 45041
$ws.Range("E780").Value = 4
$ws.Range("F780").Value = 100112024
$ws.Range("G780").Value = 'Choclo'
$ws.Range("H780").Value = 'Dulce o Americano'
$ws.Range("I780").Value = 'Primera'
$ws.Range("J780").Value = 11000
$ws.Range("K780").Value = 280
$ws.Range("L780").Value = 300
$ws.Range("M780").Value = 290
$ws.Range("N780").Value = '$/unidad'
$ws.Range("O780").Value = 'Provincia de Limarí'
$ws.Range("P780").Value = 290
$ws.Range("Q780").Value = 1
$ws.Range("R780").Value = 'Hortaliza'
